
$wb = $excel.ActiveWorkbook

# ---- Sheet 3: CypherOutput_Message (duplicate of Message) ----
$msg = $wb.Worksheets.Item("Message")
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$cypherMsg = $wb.Worksheets.Add($null, $last)
$cypherMsg.Name = "CypherOutput_Message"

$msg.Range("A1:A10").Copy()
$cypherMsg.Range("A1").PasteSpecial(-4163)

# ---- Sheet 4: StatOutput ----
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$statOutput = $wb.Worksheets.Add($null, $last)
$statOutput.Name = "StatOutput"

$statOutput.Range("A1").Value = "number_of_files"
$statOutput.Range("B1").Value = "number_of_sample"
$statOutput.Range("C1").Value = "number_of_cases"
$statOutput.Range("D1").Value = "number_of_study"

# Numeric-looking counters need to land as *text* shared strings (matching
# how this workbook's exporter wrote every other cell), not as Excel
# numbers. Stage them with a quote-prefix on a scratch cell, then copy /
# paste-values so the destination cells pick up the text type cleanly.
$statOutput.Range("Z1").Value = "'1"
$statOutput.Range("Z2").Value = "'2"

$statOutput.Range("Z1").Copy()
$statOutput.Range("A2").PasteSpecial(-4163)
$statOutput.Range("Z2").Copy()
$statOutput.Range("B2").PasteSpecial(-4163)
$statOutput.Range("C2").PasteSpecial(-4163)
$statOutput.Range("D2").PasteSpecial(-4163)

$statOutput.Range("Z1:Z2").Clear()

# ---- Sheet 5: StatOutput_Message ----
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$statMsg = $wb.Worksheets.Add($null, $last)
$statMsg.Name = "StatOutput_Message"

$msg.Range("A1:A10").Copy()
$statMsg.Range("A1").PasteSpecial(-4163)
$statMsg.Range("A1:A10").Copy()
$statMsg.Range("A11").PasteSpecial(-4163)

$statMsg.Range("A18").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN[''Vizsla'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'

$wb.Worksheets.Item("CypherOutput").Activate()
